# Add a new "GLOBE WIRELESS" sheet, cloned from the existing "FUSE" sheet,
# placed right after it, and make it the active/selected sheet (mirroring
# the workbook's new activeTab + tabSelected state from the diff).

$wb = $excel.ActiveWorkbook

$fuse = $wb.Worksheets.Item("FUSE")

# Copy FUSE to a position right after itself; Excel names the clone
# "FUSE (2)" and makes it the active sheet/tab.
$fuse.Copy($null, $fuse)

$newSheet = $wb.ActiveSheet
$newSheet.Name = "GLOBE WIRELESS"

# Match the new sheet's stored selection (activeCell/sqref = D31).
[void]$newSheet.Range("D31").Select()
